$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 565.8889
$ws.Range("I38").Value = 565.8889
$ws.Range("K38").Value = 1697.6667
$ws.Range("M38").Value = -1325.6667
$ws.Range("H58").Value = 4002.5
$ws.Range("I58").Value = 75
$ws.Range("J58").Value = 5685.7144
$ws.Range("K58").Value = 225
$ws.Range("L58").Value = 17057.1432
$ws.Range("M58").Value = -75
$ws.Range("N58").Value = -17357.1432
$ws.Range("H76").Value = 7262.1875
$ws.Range("I76").Value = 6313.7144
$ws.Range("J76").Value = 7999.8887
$ws.Range("K76").Value = 6313.7144
$ws.Range("L76").Value = 7999.8887
$ws.Range("M76").Value = -5998.7144
$ws.Range("N76").Value = -8629.8887
$ws.Range("H79").Value = 7262.1875
$ws.Range("I79").Value = 6313.7144
$ws.Range("J79").Value = 7999.8887
$ws.Range("K79").Value = 6313.7144
$ws.Range("L79").Value = 7999.8887
$ws.Range("M79").Value = -5221.7144
$ws.Range("N79").Value = -10183.8887
$ws.Range("H87").Value = 80002
$ws.Range("J87").Value = 80002
$ws.Range("L87").Value = 80002
$ws.Range("N87").Value = -82498
$ws.Range("H90").Value = 80002
$ws.Range("J90").Value = 80002
$ws.Range("L90").Value = 240006
$ws.Range("N90").Value = -252486
$ws.Range("H113").Value = 3633.7778
$ws.Range("I113").Value = 3643.4285
$ws.Range("K113").Value = 3643.4285
$ws.Range("M113").Value = -389.4285
$ws.Range("H135").Value = 2374.75
$ws.Range("J135").Value = 3249.5
$ws.Range("L135").Value = 29245.5
$ws.Range("N135").Value = -34315.5
$ws.Range("H141").Value = 10796.75
$ws.Range("I141").Value = 10796.75
$ws.Range("K141").Value = 32390.25
$ws.Range("M141").Value = -27210.25

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16115.532
$ws.Range("I32").Value = 7207.516
$ws.Range("J32").Value = 25023.549
$ws.Range("K32").Value = 7207.516
$ws.Range("L32").Value = 25023.549
$ws.Range("M32").Value = -6920.516
$ws.Range("N32").Value = -25597.549
$ws.Range("H44").Value = 30000
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H45").Value = 2013
$ws.Range("I45").Value = 2017.3334
$ws.Range("K45").Value = 2017.3334
$ws.Range("M45").Value = -1640.3334
$ws.Range("H46").Value = 26666.666
$ws.Range("I46").Value = 26666.666
$ws.Range("K46").Value = 26666.666
$ws.Range("M46").Value = -26347.666
$ws.Range("H61").Value = 1913.3846
$ws.Range("I61").Value = 1781.1666
$ws.Range("J61").Value = 3500
$ws.Range("K61").Value = 1781.1666
$ws.Range("L61").Value = 3500
$ws.Range("M61").Value = -1569.1666
$ws.Range("N61").Value = -3924
$ws.Range("H63").Value = 8222.223
$ws.Range("I63").Value = 8666.666999999999
$ws.Range("J63").Value = 8000
$ws.Range("K63").Value = 8666.666999999999
$ws.Range("L63").Value = 8000
$ws.Range("M63").Value = -7980.666999999999
$ws.Range("N63").Value = -9372
$ws.Range("H66").Value = 8222.223
$ws.Range("I66").Value = 8666.666999999999
$ws.Range("J66").Value = 8000
$ws.Range("K66").Value = 43333.335
$ws.Range("L66").Value = 40000
$ws.Range("M66").Value = -39901.335
$ws.Range("N66").Value = -46864
$ws.Range("H74").Value = 2289.9333
$ws.Range("I74").Value = 1009.5714
$ws.Range("K74").Value = 1009.5714
$ws.Range("M74").Value = -135.5714
$ws.Range("H77").Value = 2289.9333
$ws.Range("I77").Value = 1009.5714
$ws.Range("K77").Value = 5047.857
$ws.Range("M77").Value = -679.857
$ws.Range("H80").Value = 85000
$ws.Range("J80").Value = 85000
$ws.Range("L80").Value = 85000
$ws.Range("N80").Value = -86996
$ws.Range("H83").Value = 85000
$ws.Range("J83").Value = 85000
$ws.Range("L83").Value = 255000
$ws.Range("N83").Value = -264984
$ws.Range("H97").Value = 666.61536
$ws.Range("I97").Value = 659.6667
$ws.Range("J97").Value = 750
$ws.Range("K97").Value = 659.6667
$ws.Range("L97").Value = 750
$ws.Range("M97").Value = -163.6667
$ws.Range("N97").Value = -1742
$ws.Range("H136").Value = 1913.3846
$ws.Range("I136").Value = 1781.1666
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 5343.4998
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -2793.4998
$ws.Range("N136").Value = -15600

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 7500
$ws.Range("I82").Value = 7500
$ws.Range("K82").Value = 7500
$ws.Range("M82").Value = -7117
$ws.Range("H85").Value = 7500
$ws.Range("I85").Value = 7500
$ws.Range("K85").Value = 7500
$ws.Range("M85").Value = -6174

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4058.8696
$ws.Range("J31").Value = 4924.7334
$ws.Range("L31").Value = 4924.7334
$ws.Range("N31").Value = -5514.7334
$ws.Range("H34").Value = 4058.8696
$ws.Range("J34").Value = 4924.7334
$ws.Range("L34").Value = 4924.7334
$ws.Range("N34").Value = -5328.7334
$ws.Range("H58").Value = 5243.6
$ws.Range("I58").Value = 3553
$ws.Range("K58").Value = 3553
$ws.Range("M58").Value = -3350
$ws.Range("H99").Value = 10357.296
$ws.Range("I99").Value = 8245.857
$ws.Range("J99").Value = 12631.154
$ws.Range("K99").Value = 8245.857
$ws.Range("L99").Value = 12631.154
$ws.Range("M99").Value = -6747.857
$ws.Range("N99").Value = -15627.154
$ws.Range("H105").Value = 3198.077
$ws.Range("I105").Value = 1775.6
$ws.Range("J105").Value = 4087.125
$ws.Range("K105").Value = 1775.6
$ws.Range("L105").Value = 4087.125
$ws.Range("M105").Value = -28.59999999999991
$ws.Range("N105").Value = -7581.125
$ws.Range("H126").Value = 10357.296
$ws.Range("I126").Value = 8245.857
$ws.Range("J126").Value = 12631.154
$ws.Range("K126").Value = 24737.571
$ws.Range("L126").Value = 37893.462
$ws.Range("M126").Value = -22267.571
$ws.Range("N126").Value = -42833.462
$ws.Range("H136").Value = 5243.6
$ws.Range("I136").Value = 3553
$ws.Range("K136").Value = 10659
$ws.Range("M136").Value = -8109

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1679.6364
$ws.Range("J34").Value = 2235
$ws.Range("L34").Value = 6705
$ws.Range("N34").Value = -6873
$ws.Range("H39").Value = 3375.375
$ws.Range("J39").Value = 4400.6
$ws.Range("L39").Value = 13201.8
$ws.Range("N39").Value = -13789.8
$ws.Range("H55").Value = 1658
$ws.Range("J55").Value = 1821.6666
$ws.Range("L55").Value = 5464.9998
$ws.Range("N55").Value = -5818.9998
$ws.Range("H122").Value = 918.1
$ws.Range("I122").Value = 772.5
$ws.Range("J122").Value = 1136.5
$ws.Range("K122").Value = 6952.5
$ws.Range("L122").Value = 10228.5
$ws.Range("M122").Value = -4502.5
$ws.Range("N122").Value = -15128.5
$ws.Range("H136").Value = 12793.333
$ws.Range("I136").Value = 12793.333
$ws.Range("K136").Value = 38379.999
$ws.Range("M136").Value = -33279.999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6885.2
$ws.Range("I80").Value = 6106.75
$ws.Range("J80").Value = 9999
$ws.Range("K80").Value = 6106.75
$ws.Range("L80").Value = 9999
$ws.Range("M80").Value = -5108.75
$ws.Range("N80").Value = -11995
$ws.Range("H83").Value = 6885.2
$ws.Range("I83").Value = 6106.75
$ws.Range("J83").Value = 9999
$ws.Range("K83").Value = 30533.75
$ws.Range("L83").Value = 49995
$ws.Range("M83").Value = -25541.75
$ws.Range("N83").Value = -59979
$ws.Range("H126").Value = 3005.353
$ws.Range("I126").Value = 2335.3333
$ws.Range("J126").Value = 3370.818
$ws.Range("K126").Value = 7005.999899999999
$ws.Range("L126").Value = 10112.454
$ws.Range("M126").Value = -4535.999899999999
$ws.Range("N126").Value = -15052.454
$ws.Range("H132").Value = 2840.7778
$ws.Range("I132").Value = 2231.8572
$ws.Range("K132").Value = 6695.571599999999
$ws.Range("M132").Value = -4165.571599999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3107.4
$ws.Range("I7").Value = 2637.125
$ws.Range("J7").Value = 4988.5
$ws.Range("K7").Value = 2637.125
$ws.Range("L7").Value = 4988.5
$ws.Range("M7").Value = -2525.125
$ws.Range("N7").Value = -5212.5
$ws.Range("H40").Value = 2200.4
$ws.Range("I40").Value = 1858.2858
$ws.Range("K40").Value = 1858.2858
$ws.Range("M40").Value = -1722.2858
$ws.Range("H122").Value = 3534.7
$ws.Range("I122").Value = 3507.75
$ws.Range("K122").Value = 10523.25
$ws.Range("M122").Value = -8073.25
$ws.Range("H126").Value = 3107.4
$ws.Range("I126").Value = 2637.125
$ws.Range("J126").Value = 4988.5
$ws.Range("K126").Value = 7911.375
$ws.Range("L126").Value = 14965.5
$ws.Range("M126").Value = -5441.375
$ws.Range("N126").Value = -19905.5
$ws.Range("H136").Value = 2640.7144
$ws.Range("I136").Value = 2640.7144
$ws.Range("K136").Value = 7922.1432
$ws.Range("M136").Value = -5372.1432

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2475.7778
$ws.Range("J113").Value = 2877.6667
$ws.Range("L113").Value = 8633.000100000001
$ws.Range("N113").Value = -12973.0001
$ws.Range("H132").Value = 2152.5833
$ws.Range("I132").Value = 1783.2
$ws.Range("K132").Value = 5349.6
$ws.Range("M132").Value = -2819.6
